# Update column G ("K") values on the active worksheet to reflect the
# regenerated strikeout counts (K) computed from the underlying save data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 0
    6  = 2
    7  = 3
    8  = 3
    9  = 0
    10 = 1
    11 = 2
    12 = 2
    13 = 0
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    19 = 0
    20 = 1
    21 = 2
    22 = 0
    23 = 2
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 1
    29 = 3
    30 = 1
    31 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
